# Updates cryptos list: prices and 1h volume % changes for several rows,
# plus a few coin rows (Avalanche <-> WrappedliquidstakedEther2.0,
# Bittensor <-> Aave) that swapped rank positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.119.28'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +1.10%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.353.36'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.92'
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = '  +0.64%  '

$ws.Range("E8").Value = '  +5.17%  '

$ws.Range("E9").Value = '  +3.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.53'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +2.58%  '

$ws.Range("E11").Value = '  -1.95%  '

$ws.Range("E12").Value = '  -0.49%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'

$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.771.03'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +0.88%  '

$ws.Range("B14").Value = 'Avalanche'

$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.78'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.068.51'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +0.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000136'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +1.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.350.86'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("E18").Value = '  +2.45%  '

$ws.Range("E19").Value = '  +2.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '329.42'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -1.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.74'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +0.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.35'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +3.08%  '

$ws.Range("E24").Value = '  -2.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +0.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.26'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -2.32%  '

$ws.Range("E27").Value = '  -5.02%  '

$ws.Range("E28").Value = '  -0.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.36'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +0.29%  '

$ws.Range("E30").Value = '  +0.46%  '

$ws.Range("E31").Value = '  -0.42%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.37'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -0.72%  '

$ws.Range("E33").Value = '  -2.34%  '

$ws.Range("E35").Value = '  +0.86%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.17'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -0.43%  '

$ws.Range("E37").Value = '  -2.21%  '

$ws.Range("E38").Value = '  -2.39%  '

$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("B40").Value = 'Aave'

$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '141.03'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -5.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.65'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +0.91%  '

$ws.Range("B42").Value = 'Bittensor'

$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '289.37'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +0.35%  '

$ws.Range("E43").Value = '  +2.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0514'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +1.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.98'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -1.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.567'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +1.14%  '

$ws.Range("E47").Value = '  +2.13%  '

$ws.Range("E48").Value = '  -0.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.09'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("E50").Value = '  +0.62%  '

$ws.Range("E51").Value = '  +0.59%  '
